# ===========================================================
# Edit script: add Telangana row, rename Andaman & Nicobar,
# and add a new "Dataset name" column C with per-row mappings.
# ===========================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Andaman and Nicobar label (row 3, col A) ---
$ws.Range("A3").Value = "Andaman and Nicobar Islands"

# --- Make room for the new "Telangana" row by shifting rows 34-37 down to 35-38 ---
# Duplicate formatting (incl. borders/fill) of row 37 into the newly-needed row 38
# using a formats-only paste so no new style entries are created.
$ws.Range("A37:B37").Copy()
$ws.Range("A38:B38").PasteSpecial(-4122)
$ws.Rows.Item(38).RowHeight = 32.1

# Cache the values that need to move down one row (read via .Value() to invoke the getter)
$a34 = $ws.Range("A34").Value()
$b34 = $ws.Range("B34").Value()
$a35 = $ws.Range("A35").Value()
$b35 = $ws.Range("B35").Value()
$a36 = $ws.Range("A36").Value()
$b36 = $ws.Range("B36").Value()
$a37 = $ws.Range("A37").Value()
$b37 = $ws.Range("B37").Value()

# Write the cached values one row further down (bottom-up to avoid clobbering)
$ws.Range("A38").Value = $a37
$ws.Range("B38").Value = $b37
$ws.Range("A37").Value = $a36
$ws.Range("B37").Value = $b36
$ws.Range("A36").Value = $a35
$ws.Range("B36").Value = $b35
$ws.Range("A35").Value = $a34
$ws.Range("B35").Value = $b34

# --- Populate the new Telangana row (34) ---
$ws.Range("A34").Value = "Telangana`n"
$ws.Range("B34").Value = "TE`n"

# --- Header for new column C ---
$ws.Range("C2").Value = "Dataset name"

# --- Populate column C ("Dataset name") for every data row ---
$ws.Range("C3").Value = "Andaman and Nicobar Islands"
$ws.Range("C4").Value = "Andhra Pradesh"
$ws.Range("C5").Value = "Arunachal Pradesh"
$ws.Range("C6").Value = "Assam"
$ws.Range("C7").Value = "Bihar"
$ws.Range("C8").Value = "Chandigarh"
$ws.Range("C9").Value = "Chhattisgarh"
$ws.Range("C10").Value = "NOT CONSIDERED"
$ws.Range("C11").Value = "NOT CONSIDERED"
$ws.Range("C12").Value = "Delhi"
$ws.Range("C13").Value = "Goa"
$ws.Range("C14").Value = "Gujarat"
$ws.Range("C15").Value = "Haryana"
$ws.Range("C16").Value = "Himachal Pradesh"
$ws.Range("C17").Value = "Jammu and Kashmir"
$ws.Range("C18").Value = "Jharkhand"
$ws.Range("C19").Value = "Karnataka"
$ws.Range("C20").Value = "Kerala"
$ws.Range("C21").Value = "Lakshadweep"
$ws.Range("C22").Value = "Madhya Pradesh"
$ws.Range("C23").Value = "Maharashtra"
$ws.Range("C24").Value = "Manipur"
$ws.Range("C25").Value = "Meghalaya"
$ws.Range("C26").Value = "Mizoram"
$ws.Range("C27").Value = "Nagaland"
$ws.Range("C28").Value = "Odisha"
$ws.Range("C29").Value = "Puducherry"
$ws.Range("C30").Value = "Punjab"
$ws.Range("C31").Value = "Rajasthan"
$ws.Range("C32").Value = "Sikkim"
$ws.Range("C33").Value = "Tamil Nadu"
$ws.Range("C34").Value = "Telangana"
$ws.Range("C35").Value = "Tripura"
$ws.Range("C36").Value = "Uttar Pradesh"
$ws.Range("C37").Value = "Uttarakhand"
$ws.Range("C38").Value = "West Bengal"

# --- Update the active selection to match the authored state (cell C11) ---
$ws.Range("C11").Select()

